# Applies the openpyxl-authored change described by the diff:
#  - Shift the header row left by one column (A1:G1), filling in the
#    previously-empty A1 with "STUDENTID" and dropping the old H1
#    ("Percentage" now lives in G1).
#  - Populate five new data rows (2-6) with student/course records.
#  - Colour row 2 (the first record) with a solid red fill and rows 3-6
#    with a solid green fill, matching the two new cellXfs/fills added
#    to styles.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (shifted one column to the left) ---------------------
$ws.Range("A1").Value = "STUDENTID"
$ws.Range("B1").Value = "Studentname"
$ws.Range("C1").Value = "Marks"
$ws.Range("D1").Value = "courseId"
$ws.Range("E1").Value = "courseName"
$ws.Range("F1").Value = "courseDescription"
$ws.Range("G1").Value = "Percentage"

# Drop the old trailing column entirely so it no longer counts towards
# the sheet's used range / dimension.
$ws.Range("H1").Clear()

# --- Data rows ---------------------------------------------------------
$rows = @(
    @(1, "wardha", 45, 1, "Distributed Programming", "This course is related to programming", 45),
    @(5, "adnan",  90, 1, "Distributed Programming", "This course is related to programming", 90),
    @(2, "Kashif", 85, 2, "networking",               "This course is related to Networking",  85),
    @(3, "wardha", 67, 3, "Databases",                 "This course is related to Databases",   67),
    @(4, "Kashif", 88, 4, "DSA",                       "This course is related to DSA",          88)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# --- Conditional-looking static fills ----------------------------------
# Row 2 -> solid red, rows 3-6 -> solid green (RGB as BGR-packed ints,
# matching VBA/COM's Interior.Color convention: R | (G<<8) | (B<<16)).
$ws.Range("A2:G2").Interior.Color = 255      # 0x0000FF -> R=FF,G=00,B=00 (red)
$ws.Range("A3:G6").Interior.Color = 32768    # 0x008000 -> R=00,G=80,B=00 (green)
